$wb = $excel.ActiveWorkbook

# --- "Prix (2)" sheet: filter the "Prix" table down to Pays = "Belgique" ---
$wsPrix = $wb.Worksheets.Item("Prix (2)")
$loPrix = $wsPrix.ListObjects.Item("Prix")

# Apply a values-style AutoFilter on the "Pays" column (2nd column of the
# table / table column index 2) so only "Belgique" rows remain visible.
# This also hides all the non-matching data rows, same as Excel does when a
# filter is applied through the UI.
$loPrix.Range.AutoFilter(2, @("Belgique"), 7)

# The filtered rows previously carried an explicit "General" number-format
# style (s="1") on column B; drop that leftover formatting back to the
# default "Normal" style, matching a plain re-save after filtering.
$wsPrix.Range("B2:B47").Style = "Normal"

# --- "TRANSPORTEUR" sheet: drop the same leftover explicit style ---
$wsTransporteur = $wb.Worksheets.Item("TRANSPORTEUR")
$wsTransporteur.Range("B2:B16").Style = "Normal"
$wsTransporteur.Range("E2:E16").Style = "Normal"
